$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.941.22"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").Value = "1.598.97"
$ws.Range("E3").Value = "  +2.74%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'212.38"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.485"
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("D10").Value = "'18.01"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").Value = "'0.0812"
$ws.Range("E11").Value = "  +4.00%  "
$ws.Range("D12").Value = "1.823.72"
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("D13").Value = "1.599.37"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").Value = "'4.00"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "25.961.41"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").Value = "'60.18"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "0.0₃0720"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "'200.84"
$ws.Range("E20").Value = "  +8.22%  "
$ws.Range("D21").Value = "'4.21"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").Value = "'9.24"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "'5.98"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("E24").Value = "  +10.24%  "
$ws.Range("D25").Value = "'141.13"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'0.123"
$ws.Range("E27").Value = "  -5.54%  "
$ws.Range("D28").Value = "'15.12"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("E36").Value = "  +11.49%  "
$ws.Range("D37").Value = "1.127.29"
$ws.Range("E37").Value = "  +3.71%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "'0.788"
$ws.Range("E39").Value = "  +3.00%  "
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("D41").Value = "'0.488"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("D42").Value = "'0.780"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("D43").Value = "1.735.85"
$ws.Range("E43").Value = "  +2.92%  "
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").Value = "'92.93"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("E46").Value = "  +3.09%  "
$ws.Range("D47").Value = "'53.24"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "0.0₇0925"
$ws.Range("E51").Value = "  -16.78%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
